# The source diff for this revision consists entirely of XML attribute /
# namespace-declaration re-ordering (e.g. <w:pgSz w:w=".." w:h=".."/> becoming
# <w:pgSz w:h=".." w:w=".."/>, latentStyles/lsdException attributes being
# alphabetized, style element attributes being alphabetized, etc.) produced
# by the upgrade of the authoring library (Apache POI 3.15) that generated
# this fixture. Every removed line and its corresponding added line carry
# the exact same element name and the exact same set of attribute
# name/value pairs -- only the serialization order changed. There is no
# textual, structural, formatting, or style content to modify through the
# Word object model: the document's content, styles and formatting already
# match the target, so no edits are required here.
$d = $word.ActiveDocument
